# Update the "Förändrad" (Changed) date column from 2024-05-22 to 2024-05-23
# for all data rows (rows 2-28) on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45435
}
